$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2_B = New-Object 'object[,]' 1,3
$row2_B[0,0] = 11.42425601397388
$row2_B[0,1] = 12.57167033438009
$row2_B[0,2] = 5.221001886499133
$ws.Range("B2:D2").Value = $row2_B
$row2_F = New-Object 'object[,]' 1,3
$row2_F[0,0] = 24.68600013590373
$row2_F[0,1] = 29.42284862140498
$row2_F[0,2] = 14.54135836258936
$ws.Range("F2:H2").Value = $row2_F
$row2_K = New-Object 'object[,]' 1,5
$row2_K[0,0] = 7.853958302378344
$row2_K[0,1] = 11.07648354987046
$row2_K[0,2] = 14.32418650231542
$row2_K[0,3] = 19.43977965373681
$row2_K[0,4] = 22.20713456453913
$ws.Range("K2:O2").Value = $row2_K

$row3_B = New-Object 'object[,]' 1,3
$row3_B[0,0] = 11.18718396961886
$row3_B[0,1] = 12.60011006660975
$row3_B[0,2] = 5.148500111879006
$ws.Range("B3:D3").Value = $row3_B
$row3_F = New-Object 'object[,]' 1,3
$row3_F[0,0] = 24.71898460247428
$row3_F[0,1] = 29.48018292246616
$row3_F[0,2] = 14.58168645098317
$ws.Range("F3:H3").Value = $row3_F
$row3_K = New-Object 'object[,]' 1,5
$row3_K[0,0] = 7.602609276690671
$row3_K[0,1] = 11.08405700405857
$row3_K[0,2] = 14.28926818522747
$row3_K[0,3] = 19.495225239082
$row3_K[0,4] = 22.2708171717873
$ws.Range("K3:O3").Value = $row3_K

$row4_B = New-Object 'object[,]' 1,3
$row4_B[0,0] = 11.04102812975981
$row4_B[0,1] = 12.61854715414623
$row4_B[0,2] = 5.102750786591749
$ws.Range("B4:D4").Value = $row4_B
$row4_F = New-Object 'object[,]' 1,3
$row4_F[0,0] = 24.7449236454132
$row4_F[0,1] = 29.52332812276712
$row4_F[0,2] = 14.60841699478755
$ws.Range("F4:H4").Value = $row4_F
$row4_K = New-Object 'object[,]' 1,5
$row4_K[0,0] = 7.442433046835633
$row4_K[0,1] = 11.09038744588346
$row4_K[0,2] = 14.2699353471147
$row4_K[0,3] = 19.53089529311212
$row4_K[0,4] = 22.31394662262939
$ws.Range("K4:O4").Value = $row4_K

$row5_B = New-Object 'object[,]' 1,3
$row5_B[0,0] = 10.98140415981395
$row5_B[0,1] = 12.62630625269818
$row5_B[0,2] = 5.083807582411389
$ws.Range("B5:D5").Value = $row5_B
$row5_F = New-Object 'object[,]' 1,3
$row5_F[0,0] = 24.75692275821621
$row5_F[0,1] = 29.54290241891097
$row5_F[0,2] = 14.61980530642883
$ws.Range("F5:H5").Value = $row5_F
$row5_K = New-Object 'object[,]' 1,5
$row5_K[0,0] = 7.375748190972452
$row5_K[0,1] = 11.09339028974702
$row5_K[0,2] = 14.26259278436619
$row5_K[0,3] = 19.54584120718376
$row5_K[0,4] = 22.33253425626172
$ws.Range("K5:O5").Value = $row5_K

$row6_B = New-Object 'object[,]' 1,3
$row6_B[0,0] = 10.97150228333958
$row6_B[0,1] = 12.62760951271201
$row6_B[0,2] = 5.080644259478302
$ws.Range("B6:D6").Value = $row6_B
$row6_F = New-Object 'object[,]' 1,3
$row6_F[0,0] = 24.75900145040315
$row6_F[0,1] = 29.54627289221888
$row6_F[0,2] = 14.62172625633905
$ws.Range("F6:H6").Value = $row6_F
$row6_K = New-Object 'object[,]' 1,5
$row6_K[0,0] = 7.364591781870179
$row6_K[0,1] = 11.0939144818832
$row6_K[0,2] = 14.2614060771456
$row6_K[0,3] = 19.54834776624432
$row6_K[0,4] = 22.33568181665123
$ws.Range("K6:O6").Value = $row6_K

$row7_B = New-Object 'object[,]' 1,3
$row7_B[0,0] = 11.04022416667223
$row7_B[0,1] = 12.61865079978831
$row7_B[0,2] = 5.102496513032614
$ws.Range("B7:D7").Value = $row7_B
$row7_F = New-Object 'object[,]' 1,3
$row7_F[0,0] = 24.74507968663029
$row7_F[0,1] = 29.52358404870098
$row7_F[0,2] = 14.60856857530435
$ws.Range("F7:H7").Value = $row7_F
$row7_K = New-Object 'object[,]' 1,5
$row7_K[0,0] = 7.441539345996439
$row7_K[0,1] = 11.09042622928783
$row7_K[0,2] = 14.26983414617172
$row7_K[0,3] = 19.53109519702529
$row7_K[0,4] = 22.31419320545674
$ws.Range("K7:O7").Value = $row7_K

$row8_B = New-Object 'object[,]' 1,3
$row8_B[0,0] = 11.34268667626639
$row8_B[0,1] = 12.58127440107127
$row8_B[0,2] = 5.196265097045822
$ws.Range("B8:D8").Value = $row8_B
$row8_F = New-Object 'object[,]' 1,3
$row8_F[0,0] = 24.69619240628277
$row8_F[0,1] = 29.44096655404044
$row8_F[0,2] = 14.55485501301343
$ws.Range("F8:H8").Value = $row8_F
$row8_K = New-Object 'object[,]' 1,5
$row8_K[0,0] = 7.768539971895626
$row8_K[0,1] = 11.0787465236016
$row8_K[0,2] = 14.31171296758119
$row8_K[0,3] = 19.4585604581819
$row8_K[0,4] = 22.22825580269688
$ws.Range("K8:O8").Value = $row8_K

$row9_B = New-Object 'object[,]' 1,3
$row9_B[0,0] = 11.92752164338284
$row9_B[0,1] = 12.5156849084152
$row9_B[0,2] = 5.369939286217389
$ws.Range("B9:D9").Value = $row9_B
$row9_F = New-Object 'object[,]' 1,3
$row9_F[0,0] = 24.64547636400982
$row9_F[0,1] = 29.34214915970795
$row9_F[0,2] = 14.46513231861802
$ws.Range("F9:H9").Value = $row9_F
$row9_K = New-Object 'object[,]' 1,5
$row9_K[0,0] = 8.361118691727416
$row9_K[0,1] = 11.06914510689481
$row9_K[0,2] = 14.41028204409235
$row9_K[0,3] = 19.32917038636144
$row9_K[0,4] = 22.09173042808853
$ws.Range("K9:O9").Value = $row9_K

$row10_B = New-Object 'object[,]' 1,3
$row10_B[0,0] = 12.34756248906867
$row10_B[0,1] = 12.4721512134092
$row10_B[0,2] = 5.490777326150647
$ws.Range("B10:D10").Value = $row10_B
$row10_F = New-Object 'object[,]' 1,3
$row10_F[0,0] = 24.63575609069876
$row10_F[0,1] = 29.30827615582562
$row10_F[0,2] = 14.40871083220812
$ws.Range("F10:H10").Value = $row10_F
$row10_K = New-Object 'object[,]' 1,5
$row10_K[0,0] = 8.7641566720346
$row10_K[0,1] = 11.07015520332941
$row10_K[0,2] = 14.49233439795743
$row10_K[0,3] = 19.24186577333615
$row10_K[0,4] = 22.0109792311403
$ws.Range("K10:O10").Value = $row10_K

$row11_B = New-Object 'object[,]' 1,3
$row11_B[0,0] = 12.53566293304448
$row11_B[0,1] = 12.45334835907614
$row11_B[0,2] = 5.544175230982387
$ws.Range("B11:D11").Value = $row11_B
$row11_F = New-Object 'object[,]' 1,3
$row11_F[0,0] = 24.63730607025513
$row11_F[0,1] = 29.30130239757334
$row11_F[0,2] = 14.38510174910166
$ws.Range("F11:H11").Value = $row11_F
$row11_K = New-Object 'object[,]' 1,5
$row11_K[0,0] = 8.940027921644349
$row11_K[0,1] = 11.072353866996
$row11_K[0,2] = 14.53166259542863
$row11_K[0,3] = 19.20381716104755
$row11_K[0,4] = 21.97849821437208
$ws.Range("K11:O11").Value = $row11_K

$row12_B = New-Object 'object[,]' 1,3
$row12_B[0,0] = 12.60639459954856
$row12_B[0,1] = 12.44637144936119
$row12_B[0,2] = 5.564161136684841
$ws.Range("B12:D12").Value = $row12_B
$row12_F = New-Object 'object[,]' 1,3
$row12_F[0,0] = 24.63874999142309
$row12_F[0,1] = 29.29987526782716
$row12_F[0,2] = 14.37645716739137
$ws.Range("F12:H12").Value = $row12_F
$row12_K = New-Object 'object[,]' 1,5
$row12_K[0,0] = 9.005518662145116
$row12_K[0,1] = 11.07343529971913
$row12_K[0,2] = 14.54683473218772
$row12_K[0,3] = 19.18964762485944
$row12_K[0,4] = 21.96681070704014
$ws.Range("K12:O12").Value = $row12_K

$row13_B = New-Object 'object[,]' 1,3
$row13_B[0,0] = 12.59118447395792
$row13_B[0,1] = 12.44786768679931
$row13_B[0,2] = 5.559867387366432
$ws.Range("B13:D13").Value = $row13_B
$row13_F = New-Object 'object[,]' 1,3
$row13_F[0,0] = 24.63840093735704
$row13_F[0,1] = 29.30012864548114
$row13_F[0,2] = 14.37830578576286
$ws.Range("F13:H13").Value = $row13_F
$row13_K = New-Object 'object[,]' 1,5
$row13_K[0,0] = 8.991463819787251
$row13_K[0,1] = 11.07319134743368
$row13_K[0,2] = 14.54355485368512
$row13_K[0,3] = 19.19268869023488
$row13_K[0,4] = 21.96930057536848
$ws.Range("K13:O13").Value = $row13_K

$row14_B = New-Object 'object[,]' 1,3
$row14_B[0,0] = 12.54149243102081
$row14_B[0,1] = 12.45277149550055
$row14_B[0,2] = 5.545824235879773
$ws.Range("B14:D14").Value = $row14_B
$row14_F = New-Object 'object[,]' 1,3
$row14_F[0,0] = 24.63740769965293
$row14_F[0,1] = 29.30116066198772
$row14_F[0,2] = 14.38438462960112
$ws.Range("F14:H14").Value = $row14_F
$row14_K = New-Object 'object[,]' 1,5
$row14_K[0,0] = 8.945438264169665
$row14_K[0,1] = 11.07243785742356
$row14_K[0,2] = 14.53290526656409
$row14_K[0,3] = 19.20264664773589
$row14_K[0,4] = 21.97752439970413
$ws.Range("K14:O14").Value = $row14_K

$row15_B = New-Object 'object[,]' 1,3
$row15_B[0,0] = 12.51098779061276
$row15_B[0,1] = 12.45579386723137
$row15_B[0,2] = 5.53719158191026
$ws.Range("B15:D15").Value = $row15_B
$row15_F = New-Object 'object[,]' 1,3
$row15_F[0,0] = 24.63691085405711
$row15_F[0,1] = 29.30195086431683
$row15_F[0,2] = 14.38814659609254
$ws.Range("F15:H15").Value = $row15_F
$row15_K = New-Object 'object[,]' 1,5
$row15_K[0,0] = 8.917101084685218
$row15_K[0,1] = 11.07200869118488
$row15_K[0,2] = 14.52641821499041
$row15_K[0,3] = 19.20877723306828
$row15_K[0,4] = 21.98264149642613
$ws.Range("K15:O15").Value = $row15_K

$row16_B = New-Object 'object[,]' 1,3
$row16_B[0,0] = 12.33520369887637
$row16_B[0,1] = 12.47340011250087
$row16_B[0,2] = 5.487255264682241
$ws.Range("B16:D16").Value = $row16_B
$row16_F = New-Object 'object[,]' 1,3
$row16_F[0,0] = 24.63577482712965
$row16_F[0,1] = 29.30890172689434
$row16_F[0,2] = 14.41029512510827
$ws.Range("F16:H16").Value = $row16_F
$row16_K = New-Object 'object[,]' 1,5
$row16_K[0,0] = 8.752509758779125
$row16_K[0,1] = 11.07004640848347
$row16_K[0,2] = 14.48980377717797
$row16_K[0,3] = 19.2443857902469
$row16_K[0,4] = 22.01318759208531
$ws.Range("K16:O16").Value = $row16_K

$row17_B = New-Object 'object[,]' 1,3
$row17_B[0,0] = 12.2265520961272
$row17_B[0,1] = 12.48445688667707
$row17_B[0,2] = 5.456212034360445
$ws.Range("B17:D17").Value = $row17_B
$row17_F = New-Object 'object[,]' 1,3
$row17_F[0,0] = 24.63660637126891
$row17_F[0,1] = 29.31532711268135
$row17_F[0,2] = 14.42440933812069
$ws.Range("F17:H17").Value = $row17_F
$row17_K = New-Object 'object[,]' 1,5
$row17_K[0,0] = 8.649599326167408
$row17_K[0,1] = 11.06928720823174
$row17_K[0,2] = 14.46784886678062
$row17_K[0,3] = 19.26665663711804
$row17_K[0,4] = 22.03301652153522
$ws.Range("K17:O17").Value = $row17_K

$row18_B = New-Object 'object[,]' 1,3
$row18_B[0,0] = 12.16378212427408
$row18_B[0,1] = 12.49091068998166
$row18_B[0,2] = 5.438209325011503
$ws.Range("B18:D18").Value = $row18_B
$row18_F = New-Object 'object[,]' 1,3
$row18_F[0,0] = 24.63764688945655
$row18_F[0,1] = 29.31981680457212
$row18_F[0,2] = 14.43272111684464
$ws.Range("F18:H18").Value = $row18_F
$row18_K = New-Object 'object[,]' 1,5
$row18_K[0,0] = 8.589706310664349
$row18_K[0,1] = 11.0690143170624
$row18_K[0,2] = 14.45540996848417
$row18_K[0,3] = 19.27962317667666
$row18_K[0,4] = 22.04482188005134
$ws.Range("K18:O18").Value = $row18_K

$row19_B = New-Object 'object[,]' 1,3
$row19_B[0,0] = 12.14248400568675
$row19_B[0,1] = 12.49311204085715
$row19_B[0,2] = 5.432088861186936
$ws.Range("B19:D19").Value = $row19_B
$row19_F = New-Object 'object[,]' 1,3
$row19_F[0,0] = 24.63809579915119
$row19_F[0,1] = 29.32147326759582
$row19_F[0,2] = 14.43556860837045
$ws.Range("F19:H19").Value = $row19_F
$row19_K = New-Object 'object[,]' 1,5
$row19_K[0,0] = 8.569308138055376
$row19_K[0,1] = 11.06895008615982
$row19_K[0,2] = 14.45123108322518
$row19_K[0,3] = 19.28404040963335
$row19_K[0,4] = 22.04888769606495
$ws.Range("K19:O19").Value = $row19_K

$row20_B = New-Object 'object[,]' 1,3
$row20_B[0,0] = 12.238147374186
$row20_B[0,1] = 12.48327012512799
$row20_B[0,2] = 5.459531965433571
$ws.Range("B20:D20").Value = $row20_B
$row20_F = New-Object 'object[,]' 1,3
$row20_F[0,0] = 24.63645967344871
$row20_F[0,1] = 29.31456093745904
$row20_F[0,2] = 14.42288681457639
$ws.Range("F20:H20").Value = $row20_F
$row20_K = New-Object 'object[,]' 1,5
$row20_K[0,0] = 8.660627190653111
$row20_K[0,1] = 11.06935108432109
$row20_K[0,2] = 14.47016650571878
$row20_K[0,3] = 19.2642696324359
$row20_K[0,4] = 22.03086426461773
$ws.Range("K20:O20").Value = $row20_K

$row21_B = New-Object 'object[,]' 1,3
$row21_B[0,0] = 12.55610223365534
$row21_B[0,1] = 12.45132724298861
$row21_B[0,2] = 5.549955487266259
$ws.Range("B21:D21").Value = $row21_B
$row21_F = New-Object 'object[,]' 1,3
$row21_F[0,0] = 24.63767619655907
$row21_F[0,1] = 29.30082459368679
$row21_F[0,2] = 14.38259110295087
$ws.Range("F21:H21").Value = $row21_F
$row21_K = New-Object 'object[,]' 1,5
$row21_K[0,0] = 8.958987398345343
$row21_K[0,1] = 11.07265243260433
$row21_K[0,2] = 14.53602579505739
$row21_K[0,3] = 19.19971528549216
$row21_K[0,4] = 21.97509223858161
$ws.Range("K21:O21").Value = $row21_K

$row22_B = New-Object 'object[,]' 1,3
$row22_B[0,0] = 12.76097293564842
$row22_B[0,1] = 12.43128584767359
$row22_B[0,2] = 5.607680617634469
$ws.Range("B22:D22").Value = $row22_B
$row22_F = New-Object 'object[,]' 1,3
$row22_F[0,0] = 24.64346513724238
$row22_F[0,1] = 29.29892103960461
$row22_F[0,2] = 14.35797880797225
$ws.Range("F22:H22").Value = $row22_F
$row22_K = New-Object 'object[,]' 1,5
$row22_K[0,0] = 9.147512811363105
$row22_K[0,1] = 11.0762598582903
$row22_K[0,2] = 14.58069303372177
$row22_K[0,3] = 19.15891584224621
$row22_K[0,4] = 21.94221152115297
$ws.Range("K22:O22").Value = $row22_K

$row23_B = New-Object 'object[,]' 1,3
$row23_B[0,0] = 12.65191962222375
$row23_B[0,1] = 12.44190609901658
$row23_B[0,2] = 5.576999914242315
$ws.Range("B23:D23").Value = $row23_B
$row23_F = New-Object 'object[,]' 1,3
$row23_F[0,0] = 24.63991926995623
$row23_F[0,1] = 29.29928974052822
$row23_F[0,2] = 14.37095723159951
$ws.Range("F23:H23").Value = $row23_F
$row23_K = New-Object 'object[,]' 1,5
$row23_K[0,0] = 9.047495237862526
$row23_K[0,1] = 11.07420229291026
$row23_K[0,2] = 14.55670755406507
$row23_K[0,3] = 19.18056437170883
$row23_K[0,4] = 21.95943373982329
$ws.Range("K23:O23").Value = $row23_K

$row24_B = New-Object 'object[,]' 1,3
$row24_B[0,0] = 12.2329060927599
$row24_B[0,1] = 12.48380635747327
$row24_B[0,2] = 5.458031508415474
$ws.Range("B24:D24").Value = $row24_B
$row24_F = New-Object 'object[,]' 1,3
$row24_F[0,0] = 24.6365242434231
$row24_F[0,1] = 29.31490484695641
$row24_F[0,2] = 14.42357453278864
$ws.Range("F24:H24").Value = $row24_F
$row24_K = New-Object 'object[,]' 1,5
$row24_K[0,0] = 8.655643759489994
$row24_K[0,1] = 11.06932169630894
$row24_K[0,2] = 14.46911812964256
$row24_K[0,3] = 19.26534829026863
$row24_K[0,4] = 22.03183603710476
$ws.Range("K24:O24").Value = $row24_K

$row25_B = New-Object 'object[,]' 1,3
$row25_B[0,0] = 11.7706833836422
$row25_B[0,1] = 12.53260810196843
$row25_B[0,2] = 5.324105480877705
$ws.Range("B25:D25").Value = $row25_B
$row25_F = New-Object 'object[,]' 1,3
$row25_F[0,0] = 24.65435879514571
$row25_F[0,1] = 29.36209224647967
$row25_F[0,2] = 14.48773546068671
$ws.Range("F25:H25").Value = $row25_F
$row25_K = New-Object 'object[,]' 1,5
$row25_K[0,0] = 8.206315534149374
$row25_K[0,1] = 11.07032263498115
$row25_K[0,2] = 14.38189479693248
$row25_K[0,3] = 19.36280597080369
$row25_K[0,4] = 22.12523341583143
$ws.Range("K25:O25").Value = $row25_K
